$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A217").Value = "2023-12-12 08:28:12"
$ws.Range("B217").Value = 0.0018

$ws.Range("A218").Value = "2023-12-12 08:28:49"
$ws.Range("B218").Value = 0.0028

$ws.Range("A219").Value = "2023-12-12 08:29:20"
$ws.Range("B219").Value = 0.002
